# Generate Report for Handback
#
# The aa7bb647... handback for zh-cn/de-de has finished syncing with en-US,
# so the status moves from "Ready for handoff" to "Handed back: in sync with
# en-US" (this text lives on the Overview summary sheet as well as the
# per-language detail sheets), the handback timestamps for the first data
# row are refreshed, and the now-resolved "handback file is not the latest"
# error detail for the second data row is cleared. A couple of columns are
# also resized to better fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# Shared by the Overview sheet (E2:F3) and the per-language sheets' Status
# column (C2:C3).
$wsOverview.Range("E2:F3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2:C3").Value     = "Handed back: in sync with en-US"
$wsDeDe.Range("C2:C3").Value     = "Handed back: in sync with en-US"

# --- Refresh "Latest Handback DateTime" (column K). Both data rows share
#     the same timestamp string, so both get the refreshed value. ---
$wsZhCn.Range("K2:K3").Value = "2016-10-25 03:25:39"
$wsDeDe.Range("K2:K3").Value = "2016-10-25 03:25:55"

# --- Clear the stale "handback file is not the latest" error detail
#     (column P, Error Detail) now that the handback is in sync ---
$wsZhCn.Range("P3").Value = ""
$wsDeDe.Range("P3").Value = ""

# --- Column width adjustments to fit the new status text / narrower errors ---
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 29.166666666666668
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 29.166666666666668

$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 29.166666666666668
$wsZhCn.Range("P1").EntireColumn.ColumnWidth = 12.833333333333334

$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 29.166666666666668
$wsDeDe.Range("P1").EntireColumn.ColumnWidth = 12.833333333333334
